$wb = $excel.ActiveWorkbook

# --- Rename sheets: "Goals" -> "Final Products", "Scopes" -> "Phases"  ---
# (Tab rename only - the "Goals"/"Final Products" sheet's own header row is
# untouched: its A1 still reads "Project Name", referencing its parent level.)
$wb.Worksheets.Item("Goals").Name = "Final Products"
$phases = $wb.Worksheets.Item("Scopes")
$phases.Name = "Phases"

# The "Phases" (formerly "Scopes") sheet's A1 is a back-reference to its
# parent level's name, so it is reworded to match the "Final Products" rename.
$phases.Range("A1").Value = "Final Product Description"

# --- Build the new "Work Packages" sheet by copying the existing
#     "Deliverables" layout (same column widths/styles) before it is
#     trimmed, so the new sheet keeps the "Assignee" column. ---
$deliverables = $wb.Worksheets.Item("Deliverables")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$deliverables.Copy([System.Reflection.Missing]::Value, $lastSheet)

$workPackages = $wb.Worksheets.Item($wb.Worksheets.Count)
$workPackages.Name = "Work Packages"

# Work Packages keeps: Deliverable Description | Description | Assignee | Budget | Status
# -> drop the "Owner" column (D) and set the back-reference header in A1.
$workPackages.Range("D1").EntireColumn.Delete()
$workPackages.Range("A1").Value = "Deliverable Description"

# --- Trim the original "Deliverables" sheet: the "Assignee" column moves
#     to the new Work Packages sheet, so drop it here, and reword the
#     back-reference header in A1 to match the "Phases" rename. ---
$deliverables.Range("C1").EntireColumn.Delete()
$deliverables.Range("A1").Value = "Phase Description"
